$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '36.673.81'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.963.25'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.52'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.76'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +7.83%  '
$ws.Range('E9').Value = '  +5.24%  '
$ws.Range('E10').Value = '  -4.72%  '
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.35'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.61%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.839'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.00%  '
$ws.Range('E14').Value = '  +3.82%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.251.69'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('E16').Value = '  +4.00%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.958.88'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '36.603.56'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.02'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '230.55'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.05%  '
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.47'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +6.13%  '
$ws.Range('E25').Value = '  +4.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.143'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +8.59%  '
$ws.Range('E27').Value = '  +1.88%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '160.97'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.45'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.89%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.23'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +12.60%  '
$ws.Range('E31').Value = '  +2.23%  '
$ws.Range('E32').Value = '  +6.17%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0619'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  +8.15%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.59'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +21.74%  '
$ws.Range('E36').Value = '  +6.58%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.60'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.63%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0989'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.25%  '
$ws.Range('E41').Value = '  +1.16%  '
$ws.Range('E42').Value = '  +3.33%  '
$ws.Range('E43').Value = '  +1.83%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.28'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +5.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.367.96'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.40%  '
$ws.Range('E46').Value = '  +2.67%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '88.83'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.78%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.17'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.84'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '44.40'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.91%  '
$ws.Range('E51').Value = '  +6.80%  '
